# Apply scheduled-runner updates to Malboro_Profits workbook
# Writes updated market-price/profit figures into the affected leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3219.3794
$ws.Range("I98").Value = 3274.56
$ws.Range("K98").Value = 3274.56
$ws.Range("M98").Value = -1776.56

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 8105.75
$ws.Range("I113").Value = 8586.929
$ws.Range("J113").Value = 4737.5
$ws.Range("K113").Value = 8586.929
$ws.Range("L113").Value = 4737.5
$ws.Range("M113").Value = -5332.929
$ws.Range("N113").Value = -11245.5

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3219.3794
$ws.Range("I122").Value = 3274.56
$ws.Range("K122").Value = 9823.68
$ws.Range("M122").Value = -7373.68

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 247999.33
$ws.Range("I137").Value = 9999
$ws.Range("J137").Value = 277749.38
$ws.Range("K137").Value = 29997
$ws.Range("L137").Value = 833248.14
$ws.Range("M137").Value = -27447
$ws.Range("N137").Value = -838348.14

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2992.6191
$ws.Range("I45").Value = 3140.5
$ws.Range("J45").Value = 2519.4
$ws.Range("K45").Value = 3140.5
$ws.Range("L45").Value = 2519.4
$ws.Range("M45").Value = -2763.5
$ws.Range("N45").Value = -3273.4

# ARM row 93
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H93").Value = 14999
$ws.Range("I93").Value = 14999
$ws.Range("K93").Value = 14999
$ws.Range("M93").Value = -12503

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3163.3333
$ws.Range("I86").Value = 3368.1667
$ws.Range("J86").Value = 2548.8333
$ws.Range("K86").Value = 3368.1667
$ws.Range("L86").Value = 2548.8333
$ws.Range("M86").Value = -2245.1667
$ws.Range("N86").Value = -4794.8333

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3163.3333
$ws.Range("I89").Value = 3368.1667
$ws.Range("J89").Value = 2548.8333
$ws.Range("K89").Value = 16840.8335
$ws.Range("L89").Value = 12744.1665
$ws.Range("M89").Value = -11224.8335
$ws.Range("N89").Value = -23976.1665

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3043
$ws.Range("I22").Value = 4175
$ws.Range("J22").Value = 2665.6667
$ws.Range("K22").Value = 4175
$ws.Range("L22").Value = 2665.6667
$ws.Range("M22").Value = -3825
$ws.Range("N22").Value = -3365.6667

# CRP row 109
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 18000
$ws.Range("J109").Value = 18000
$ws.Range("L109").Value = 18000
$ws.Range("N109").Value = -20080

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 31255298
$ws.Range("I134").Value = 1957.84
$ws.Range("K134").Value = 5873.52
$ws.Range("M134").Value = -3338.52

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7810254
$ws.Range("I4").Value = 12076535
$ws.Range("J4").Value = 1715566.9
$ws.Range("K4").Value = 36229605
$ws.Range("L4").Value = 5146700.699999999
$ws.Range("M4").Value = -36229493
$ws.Range("N4").Value = -5146924.699999999

# CUL row 29
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1283.3
$ws.Range("I29").Value = 1520.2858
$ws.Range("J29").Value = 730.3333
$ws.Range("K29").Value = 4560.857400000001
$ws.Range("L29").Value = 2190.9999
$ws.Range("M29").Value = -4283.857400000001
$ws.Range("N29").Value = -2744.9999

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3138.875
$ws.Range("I34").Value = 901.6923
$ws.Range("J34").Value = 12833.333
$ws.Range("K34").Value = 2705.0769
$ws.Range("L34").Value = 38499.999
$ws.Range("M34").Value = -2621.0769
$ws.Range("N34").Value = -38667.999

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 200129000
$ws.Range("J37").Value = 200129000
$ws.Range("L37").Value = 600387000
$ws.Range("N37").Value = -600387224

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 17859440
$ws.Range("J55").Value = 3250
$ws.Range("L55").Value = 9750
$ws.Range("N55").Value = -10104

# CUL row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 846.5
$ws.Range("I99").Value = 697.5
$ws.Range("J99").Value = 896.1667
$ws.Range("K99").Value = 2092.5
$ws.Range("L99").Value = 2688.5001
$ws.Range("M99").Value = 153.5
$ws.Range("N99").Value = -7180.5001

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3591.2856
$ws.Range("I122").Value = 3080.75
$ws.Range("J122").Value = 3711.4119
$ws.Range("K122").Value = 27726.75
$ws.Range("L122").Value = 33402.7071
$ws.Range("M122").Value = -25276.75
$ws.Range("N122").Value = -38302.7071

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2706.8
$ws.Range("I129").Value = 1997.5
$ws.Range("K129").Value = 5992.5
$ws.Range("M129").Value = -992.5

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1478.29
$ws.Range("I131").Value = 1161.4286
$ws.Range("K131").Value = 3484.2858
$ws.Range("M131").Value = 1555.7142

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5401.161
$ws.Range("I102").Value = 6280.304
$ws.Range("J102").Value = 2873.625
$ws.Range("K102").Value = 6280.304
$ws.Range("L102").Value = 2873.625
$ws.Range("M102").Value = -4658.304
$ws.Range("N102").Value = -6117.625

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3469.5386
$ws.Range("I122").Value = 3060.4
$ws.Range("J122").Value = 4833.3335
$ws.Range("K122").Value = 9181.200000000001
$ws.Range("L122").Value = 14500.0005
$ws.Range("M122").Value = -6731.200000000001
$ws.Range("N122").Value = -19400.0005

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10049.889
$ws.Range("I7").Value = 14699.75
$ws.Range("J7").Value = 6330
$ws.Range("K7").Value = 14699.75
$ws.Range("L7").Value = 6330
$ws.Range("M7").Value = -14587.75
$ws.Range("N7").Value = -6554

# LTW row 13
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1442.4286
$ws.Range("I13").Value = 1850
$ws.Range("J13").Value = 899
$ws.Range("K13").Value = 1850
$ws.Range("L13").Value = 899
$ws.Range("M13").Value = -1710
$ws.Range("N13").Value = -1179

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1315.0667
$ws.Range("I55").Value = 1169.2
$ws.Range("J55").Value = 1460.9333
$ws.Range("K55").Value = 1169.2
$ws.Range("L55").Value = 1460.9333
$ws.Range("M55").Value = -996.2
$ws.Range("N55").Value = -1806.9333

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10575.375
$ws.Range("J122").Value = 4934
$ws.Range("L122").Value = 14802
$ws.Range("N122").Value = -19702

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 10049.889
$ws.Range("I126").Value = 14699.75
$ws.Range("J126").Value = 6330
$ws.Range("K126").Value = 44099.25
$ws.Range("L126").Value = 18990
$ws.Range("M126").Value = -41629.25
$ws.Range("N126").Value = -23930

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4201.4443
$ws.Range("I81").Value = 4645
$ws.Range("K81").Value = 9290
$ws.Range("M81").Value = -8229

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4201.4443
$ws.Range("I84").Value = 4645
$ws.Range("K84").Value = 46450
$ws.Range("M84").Value = -41146

